$wb = $excel.ActiveWorkbook

# --- Sheet references (by position, since sheet 1's name is changing) ---
$wsRuptura = $wb.Worksheets.Item(1)   # 01_BD_Ruptura_Faltaproduto -> 01_BD_Ruptura
$wsEstoque = $wb.Worksheets.Item(2)   # 02_BD_Estoque
$wsVendas  = $wb.Worksheets.Item(3)   # 03_BD_Vendas

# --- Rename the first sheet (also updates the _FilterDatabase defined name) ---
$wsRuptura.Name = "01_BD_Ruptura"

# --- Row heights: rows 1-41 get an explicit 13.2pt height on 02_BD_Estoque and 03_BD_Vendas
#     (rows 42+ already have this height in the source file) ---
$wsEstoque.Range("1:41").RowHeight = 13.2
$wsVendas.Range("1:41").RowHeight = 13.2

# --- Move the active/selected tab from 03_BD_Vendas to 01_BD_Ruptura ---
$wsRuptura.Activate()
